$wb = $excel.ActiveWorkbook

# Update the B5 value (train.iteration) on the "config" sheet from 28 to 2
$configSheet = $wb.Worksheets.Item("config")
$configSheet.Range("B5").Value = 2

# Reset the selection on the "train" sheet back to A1 (was A7)
$trainSheet = $wb.Worksheets.Item("train")
$trainSheet.Activate()
$trainSheet.Range("A1").Select()
